# Generate Report for Archive
#
# The localization status for this handoff moved from "Ready for handoff"
# to "In Translation". Update the Status cells on all three sheets, and
# narrow the Status-related columns that Excel auto-fits to the new
# (shorter) text.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Overview" ---------------------------------------------------
# Columns E ("zh-cn") and F ("de-de") hold the per-language status text.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
# Narrow both columns to match the new, shorter status text.
$wsOverview.Range("E1:F1").ColumnWidth = 12.5

# --- Sheet 2: "zh-cn" --------------------------------------------------------
# Column C ("Status") holds the status text.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C1").ColumnWidth = 12.5

# --- Sheet 3: "de-de" --------------------------------------------------------
# Column C ("Status") holds the status text.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C1").ColumnWidth = 12.5
